$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("工作表1")

# Remove the Author column entirely (shifts everything left by one)
$ws.Columns("A:A").Delete()

# Remove the redundant duplicate columns (old Y, Zr, Nb, Cs, Ba, La) - 6 columns, now at J:O
$ws.Columns("J:O").Delete()

# Update the header row
$ws.Range("A1").Value = "Label"
$ws.Range("B1").Value = "A"
$ws.Range("C1").Value = "B"
$ws.Range("D1").Value = "C"
$ws.Range("E1").Value = "D"
$ws.Range("F1").Value = "E"
$ws.Range("G1").Value = "F"
$ws.Range("H1").Value = "G"
$ws.Range("I1").Value = "H"

# Update the group labels in column A (rows 2-9)
$ws.Range("A2").Value = "Group1"
$ws.Range("A3").Value = "Group1"
$ws.Range("A4").Value = "Group1"
$ws.Range("A5").Value = "Group2"
$ws.Range("A6").Value = "Group2"
$ws.Range("A7").Value = "Group3"
$ws.Range("A8").Value = "Group3"
$ws.Range("A9").Value = "Group3"

# Row 1 height and the new view/selection state
$ws.Rows(1).RowHeight = 18
$ws.Range("J1:O1048576").Select()
